# Filter - Study - Test Suit
# Re-order / relabel the "startup" sheet rows:
#   - Remove the "CasesTab" row content, replace with a "ParticipantsTab" row
#     (using the participant-focused Cypher query).
#   - "SamplesTab" row moves up to row 3 (same content, different query text
#     that was previously on CasesTab's query column reused appropriately).
#   - "FilesTab" row moves up to row 4.
#   - Selection moves from B12 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$qParticipants = "MATCH (s:study)<--(p:participant)`nWHERE s.study_name in [`"Clonal evolution during metastatic spread in high-rish neuroblastoma`"]`nOPTIONAL MATCH (p)<--(samp:sample)`nWITH p, s, collect(distinct samp.sample_id) as samp`nRETURN   `n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`n coalesce(p.gender,'') as ``Gender``,`n coalesce(apoc.text.join(samp, ','), '') as ``Samples```n ORDER By p.participant_id LIMIT 100"

$qSamples = "MATCH (s:study)<--(p:participant)<--(samp:sample)`nWHERE s.study_name in [`"Clonal evolution during metastatic spread in high-rish neuroblastoma`"]`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID``,`n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(samp.sample_tumor_status,'') as ``Tumor``,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"

$qFiles = "MATCH (s:study)<--(p:participant)`nWHERE s.study_name in [`"Clonal evolution during metastatic spread in high-rish neuroblastoma`"]`nOPTIONAL MATCH (p)<--(samp:sample)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nWITH DISTINCT p,s,samp,f,diag`nRETURN `n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```n   ORDER By f.file_name LIMIT 100"

$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("B2").Value = $qParticipants

$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $qSamples

$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $qFiles

$ws.Range("A2").Select()
